$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "code_type"
$ws.Range("B1").Value = "code"
$ws.Range("C1").Value = "term"
$ws.Range("D1").Value = "category"

$ws.Range("J8").Select()
